# Update the "Jazdy" trip log with corrected/additional rows.
# The odometer readings must reflect the final state after each run
# (fix: "odometer final state = state is final after run"), and two
# new trips (13/14 Nov, and the 28 Nov service run) are appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jazdy")

$data = @(
    @("2025-11-03", "Vrbové",    "07:00", "Bratislava", "08:00", 125746, 92.5),
    @("2025-11-03", "Bratislava","15:30", "Vrbové",     "16:30", 125839, 92.5),
    @("2025-11-04", "Vrbové",    "06:30", "Trnava",     "07:30", 125886, 47.9),
    @("2025-11-04", "Trnava",    "14:00", "Vrbové",     "15:00", 125934, 47.9),
    @("2025-11-05", "Vrbové",    "06:45", "Žilina",     "07:45", 126060, 126.1),
    @("2025-11-05", "Žilina",    "17:00", "Vrbové",     "18:00", 126187, 126.1),
    @("2025-11-07", "Vrbové",    "06:00", "Trenčín",    "07:00", 126242, 55.9),
    @("2025-11-07", "Trenčín",   "14:30", "Vrbové",     "15:30", 126298, 55.9),
    @("2025-11-11", "Vrbové",    "06:15", "Bratislava", "07:15", 126391, 92.5),
    @("2025-11-11", "Bratislava","14:45", "Vrbové",     "15:45", 126483, 92.5),
    @("2025-11-12", "Vrbové",    "07:00", "Trnava",     "08:00", 126531, 47.9),
    @("2025-11-12", "Trnava",    "15:00", "Vrbové",     "16:00", 126579, 47.9),
    @("2025-11-13", "Vrbové",    "06:30", "Žilina",     "07:30", 126705, 126.1),
    @("2025-11-13", "Žilina",    "17:00", "Vrbové",     "18:00", 126831, 126.1),
    @("2025-11-14", "Vrbové",    "07:00", "Banská Bystrica", "08:00", 127034, 202.7),
    @("2025-11-14", "Banská Bystrica", "18:00", "Vrbové", "19:00", 127237, 202.7),
    @("2025-11-28", "Vrbové",    "14:00", "Servisná Jazda (doladenie)", "15:00", 127240, 2.9),
    @("2025-11-28", "Servisná Jazda (doladenie)", "15:00", "Vrbové", "16:00", 127242, 2.9)
)

$row = 2
foreach ($r in $data) {
    # Column A holds a plain-text date string (e.g. "2025-11-03"), not a
    # native Excel date serial. Re-assigning a cell's Value re-parses the
    # string, so force a text number format first to stop Excel turning it
    # into a date serial - even when the string itself is unchanged.
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
    $row++
}
